$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing values per diff
$ws.Range("B27").Value = 2100
$ws.Range("B28").Value = 2744
$ws.Range("C28").Value = 5794

# Add new row 29
$ws.Range("A29").Value = 43857
$ws.Range("A29").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("B29").Value = 4515
$ws.Range("C29").Value = 6973

# Update selection to match diff (activeCell C29)
$ws.Range("C29").Select()
